$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), R (Origen), S (Precio $/Kg).
# Rows 3-19 have their values permuted per the source diff.
$rows = @(
    @{ Row = 3; D = 44974; M = 130; N = 7000; O = 7500; P = 7269; R = "Provincia de Curicó"; S = 3634 }
    @{ Row = 4; D = 44585; M = 160; N = 6500; O = 7000; P = 6750; R = "Provincia de Curicó"; S = 3375 }
    @{ Row = 5; D = 44589; M = 60; N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 }
    @{ Row = 6; D = 44628; M = 40; N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
    @{ Row = 7; D = 44587; M = 165; N = 6500; O = 7000; P = 6742; R = "Provincia de Linares"; S = 3371 }
    @{ Row = 8; D = 44959; M = 40; N = 7000; O = 7000; P = 7000; R = "Provincia de Curicó"; S = 3500 }
    @{ Row = 9; D = 44960; M = 40; N = 7000; O = 7000; P = 7000; R = "Provincia de Curicó"; S = 3500 }
    @{ Row = 10; D = 45001; M = 66; N = 7500; O = 8000; P = 7773; R = "Provincia de Curicó"; S = 3886 }
    @{ Row = 11; D = 44209; M = 58; N = 6000; O = 6000; P = 6000; R = "Provincia de Curicó"; S = 3000 }
    @{ Row = 12; D = 44586; M = 80; N = 7000; O = 7000; P = 7000; R = "Provincia de Curicó"; S = 3500 }
    @{ Row = 13; D = 44606; M = 45; N = 7000; O = 7000; P = 7000; R = "Provincia de Linares"; S = 3500 }
    @{ Row = 14; D = 44582; M = 150; N = 6000; O = 6500; P = 6233; R = "Provincia de Curicó"; S = 3116 }
    @{ Row = 15; D = 44588; M = 160; N = 6500; O = 7000; P = 6750; R = "Provincia de Curicó"; S = 3375 }
    @{ Row = 16; D = 44614; M = 45; N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
    @{ Row = 17; D = 44627; M = 45; N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
    @{ Row = 18; D = 44592; M = 30; N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 }
    @{ Row = 19; D = 44214; M = 48; N = 6000; O = 6000; P = 6000; R = "Provincia de Linares"; S = 3000 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("D$r").Value = $item.D
    $ws.Range("M$r").Value = $item.M
    $ws.Range("N$r").Value = $item.N
    $ws.Range("O$r").Value = $item.O
    $ws.Range("P$r").Value = $item.P
    $ws.Range("R$r").Value = $item.R
    $ws.Range("S$r").Value = $item.S
}
